# Fruta / hortaliza, semanal
# Insert two new weekly rows (186, 187) before the existing data that was
# previously at rows 186-188 (which shift down to 188-190).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 186, pushing old rows 186-188 down to 188-190
$ws.Rows.Item(186).Resize(2).Insert()

# New row 186 data
$ws.Range("A186").Value = 4
$ws.Range("B186").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C186").Value = "Los Lagos"
$ws.Range("D186").Value = 44656
$ws.Range("D186").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E186").Value = 10
$ws.Range("F186").Value = "Fruta"
$ws.Range("G186").Value = 100101
$ws.Range("H186").Value = "Berries"
$ws.Range("I186").Value = 100101007
$ws.Range("J186").Value = "Kiwi"
$ws.Range("K186").Value = "Hayward"
$ws.Range("L186").Value = "Primera"
$ws.Range("M186").Value = 600
$ws.Range("N186").Value = 17000
$ws.Range("O186").Value = 18000
$ws.Range("P186").Value = 17500
$ws.Range("Q186").Value = "`$/caja 15 kilos"
$ws.Range("R186").Value = "Región de O'Higgins"
$ws.Range("S186").Value = 1167
$ws.Range("T186").Value = 15

# New row 187 data
$ws.Range("A187").Value = 4
$ws.Range("B187").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C187").Value = "Los Lagos"
$ws.Range("D187").Value = 44656
$ws.Range("D187").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E187").Value = 10
$ws.Range("F187").Value = "Fruta"
$ws.Range("G187").Value = 100101
$ws.Range("H187").Value = "Berries"
$ws.Range("I187").Value = 100101007
$ws.Range("J187").Value = "Kiwi"
$ws.Range("K187").Value = "Hayward"
$ws.Range("L187").Value = "Segunda"
$ws.Range("M187").Value = 300
$ws.Range("N187").Value = 15000
$ws.Range("O187").Value = 15000
$ws.Range("P187").Value = 15000
$ws.Range("Q187").Value = "`$/caja 15 kilos"
$ws.Range("R187").Value = "Región de O'Higgins"
$ws.Range("S187").Value = 1000
$ws.Range("T187").Value = 15
